$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.275.91"
$ws.Range("E2").Value = "  +4.56%  "
$ws.Range("D3").Value = "3.255.14"
$ws.Range("E3").Value = "  +2.37%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.96"
$ws.Range("E5").Value = "  +2.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.31"
$ws.Range("E6").Value = "  +5.38%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.601"
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("D9").Value = "3.253.54"
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("E10").Value = "  +3.83%  "
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("E12").Value = "  +4.26%  "
$ws.Range("D13").Value = "3.828.17"
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.20"
$ws.Range("E15").Value = "  +2.73%  "
$ws.Range("D16").Value = "67.233.43"
$ws.Range("E16").Value = "  +4.47%  "
$ws.Range("E17").Value = "  +2.74%  "
$ws.Range("D18").Value = "3.258.70"
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("E20").Value = "  +2.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.63"
$ws.Range("E21").Value = "  +6.05%  "
$ws.Range("E22").Value = "  +6.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.30"
$ws.Range("E24").Value = "  +3.15%  "
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("D26").Value = "3.397.44"
$ws.Range("E26").Value = "  +2.56%  "
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.89"
$ws.Range("E28").Value = "  +3.46%  "
$ws.Range("E29").Value = "  +1.17%  "
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("E31").Value = "  +3.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.63"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("E33").Value = "  +2.35%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  +5.73%  "
$ws.Range("E36").Value = "  +2.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "164.18"
$ws.Range("E37").Value = "  +6.04%  "
$ws.Range("E38").Value = "  +4.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.855"
$ws.Range("E39").Value = "  +4.32%  "
$ws.Range("E40").Value = "  +9.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.04"
$ws.Range("E41").Value = "  +4.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.62"
$ws.Range("E42").Value = "  +2.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.60"
$ws.Range("E43").Value = "  +9.52%  "
$ws.Range("D44").Value = "2.771.60"
$ws.Range("E44").Value = "  +5.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.41"
$ws.Range("E45").Value = "  +5.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.79"
$ws.Range("E46").Value = "  +8.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "344.94"
$ws.Range("E47").Value = "  +6.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.48"
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("E49").Value = "  +2.62%  "
$ws.Range("E50").Value = "  +3.50%  "
$ws.Range("E51").Value = "  +0.87%  "
